$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix B18 / B19: the shared "=A18"/"=A19" formula is replaced by the
# literal species text (same text as before: "Mullus surmuletus" becomes
# "Mullus sp." is NOT what happens -- the row's Species now reads
# "Mullus sp." directly instead of the array-formula mirror of A18/A19).
$ws.Range("B18").Value2 = "Mullus sp."
$ws.Range("B19").Value2 = "Mullus sp."

# --- Append the 14 new data rows (118-131) describing 7 new SubSpecies,
# one F/M pair per SubSpecies.
$newRows = @(
    @("Trisopterus minutus", "Trisopterus minutus", "F"),
    @("Trisopterus minutus", "Trisopterus minutus", "M"),
    @("Trachurus trachurus", "Trachurus sp.", "F"),
    @("Trachurus trachurus", "Trachurus sp.", "M"),
    @("Trachurus mediterraneus", "Trachurus sp.", "F"),
    @("Trachurus mediterraneus", "Trachurus sp.", "M"),
    @("Eutrigla gurnardus", "Eutrigla gurnardus", "F"),
    @("Eutrigla gurnardus", "Eutrigla gurnardus", "M"),
    @("Sparus aurata", "Sparus aurata", "F"),
    @("Sparus aurata", "Sparus aurata", "M"),
    @("Boops boops", "Boops boops", "F"),
    @("Boops boops", "Boops boops", "M"),
    @("Mullus barbatus", "Mullus sp.", "F"),
    @("Mullus barbatus", "Mullus sp.", "M")
)

$startRow = 118
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value2 = $data[0]
    $ws.Cells.Item($r, 2).Value2 = $data[1]
    $ws.Cells.Item($r, 3).Value2 = $data[2]
    $ws.Cells.Item($r, 4).Value2 = "NA"
}
